$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row text updates (row 1)
$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("F1").Value = "Birthdate"
$ws.Range("G1").Value = "Ministry Course Code and Level"
$ws.Range("I1").Value = "Final Percent"
$ws.Range("K1").Value = "Credits"

# Update the selection to the header row range
$ws.Range("A1:K1").Select()
